# aggiornamento a 9/09 compreso
# Append 8 new daily rows (2021-09-02 .. 2021-09-09, serials 44441-44448)
# after the last existing row (366), copying column A's date format/border
# style from the previous row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44441, 0, 1, 21.81025081788441),
    @(44442, 0, 1, 21.81025081788441),
    @(44443, 0, 1, 21.81025081788441),
    @(44444, 0, 1, 21.81025081788441),
    @(44445, 0, 0, 0),
    @(44446, 0, 0, 0),
    @(44447, 0, 0, 0),
    @(44448, 0, 0, 0)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    # Carry the formatted (bordered, centered, date-formatted) style of
    # column A down onto the new row before writing values into it.
    $ws.Range("A" + $lastRow).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

$excel.CutCopyMode = 0
